# scripts/core_classic_mapping.xlsx — fix bug in mapping to Cognite core types
# (namespace fix: "cdf_cdm:SourceSystem" -> "cdf_cdm:CogniteSourceSystem",
#  and "Unit" -> "cdf_cdm:CogniteUnit") + add the CogniteUnit row to Views.

$wb = $excel.ActiveWorkbook

# --- Properties sheet: fix the three ClassicEvent/ClassicAsset/ClassicFile
#     "source" rows and the ClassicTimeSeries "unitExternalId" row that were
#     pointing at the wrong / non-existent core view names.
$wsProps = $wb.Worksheets.Item("Properties")
$wsProps.Range("F7").Value  = "cdf_cdm:CogniteSourceSystem(version=v1)"
$wsProps.Range("F14").Value = "cdf_cdm:CogniteSourceSystem(version=v1)"
$wsProps.Range("F21").Value = "cdf_cdm:CogniteSourceSystem(version=v1)"
$wsProps.Range("F33").Value = "cdf_cdm:CogniteUnit(version=v1)"

# --- Views sheet: replace the old (wrong) "cdf_cdm:SourceSystem(version=v1)"
#     view entry with the correct "cdf_cdm:CogniteSourceSystem(version=v1)"
#     one, and add a new row describing the CogniteUnit view.
$wsViews = $wb.Worksheets.Item("Views")
$wsViews.Range("A8").Value = "cdf_cdm:CogniteSourceSystem(version=v1)"
$wsViews.Range("C8").Value = "The CogniteSourceSystem core concept is used to standardize the way source system is stored."
$wsViews.Range("D8").Value = "cdf_cdm:CogniteDescribable(version=v1)"
$wsViews.Range("F8").Value = $true

$wsViews.Range("A9").Value = "cdf_cdm:CogniteUnit(version=v1)"
$wsViews.Range("C9").Value = "Represents a single unit of measurement"
$wsViews.Range("D9").Value = "CogniteDescribable"
$wsViews.Range("F9").Value = $true

# --- Cosmetic view-state refresh left by the authoring session ---
$wsProps.Range("F35").Select() | Out-Null
$wsViews.Range("C13").Select() | Out-Null
